$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value2 = '46.910.18'
$ws.Range('E2').Value2 = '  +0.49%  '
$ws.Range('D3').Value2 = '2.473.67'
$ws.Range('E3').Value2 = '  +0.19%  '
$ws.Range('D5').Value2 = '319.04'
$ws.Range('E5').Value2 = '  -1.33%  '
$ws.Range('D6').Value2 = '107.35'
$ws.Range('E6').Value2 = '  +1.97%  '
$ws.Range('D7').Value2 = '0.519'
$ws.Range('E7').Value2 = '  -0.63%  '
$ws.Range('E8').Value2 = '  +0.03%  '
$ws.Range('E9').Value2 = '  -1.48%  '
$ws.Range('D10').Value2 = '38.33'
$ws.Range('E10').Value2 = '  +6.14%  '
$ws.Range('D11').Value2 = '0.0804'
$ws.Range('E11').Value2 = '  -1.39%  '
$ws.Range('E12').Value2 = '  +0.14%  '
$ws.Range('D13').Value2 = '18.03'
$ws.Range('E13').Value2 = '  -1.42%  '
$ws.Range('E14').Value2 = '  -0.10%  '
$ws.Range('D15').Value2 = '2.860.42'
$ws.Range('E15').Value2 = '  +0.21%  '
$ws.Range('D16').Value2 = '2.483.28'
$ws.Range('E16').Value2 = '  -0.52%  '
$ws.Range('D17').Value2 = '0.839'
$ws.Range('E17').Value2 = '  -0.53%  '
$ws.Range('D18').Value2 = '46.840.41'
$ws.Range('E18').Value2 = '  +0.65%  '
$ws.Range('D19').Value2 = '12.63'
$ws.Range('E19').Value2 = '  -0.01%  '
$ws.Range('E20').Value2 = '  +1.47%  '
$ws.Range('D21').Value2 = '2.74'
$ws.Range('E21').Value2 = '  +15.17%  '
$ws.Range('D22').Value2 = '0.0₃0925'
$ws.Range('E22').Value2 = '  -1.21%  '
$ws.Range('D23').Value2 = '70.24'
$ws.Range('E23').Value2 = '  -0.37%  '
$ws.Range('D24').Value2 = '243.92'
$ws.Range('E24').Value2 = '  -2.05%  '
$ws.Range('D25').Value2 = '2.53'
$ws.Range('E25').Value2 = '  -0.51%  '
$ws.Range('D26').Value2 = '0.999'
$ws.Range('E26').Value2 = '  -0.06%  '
$ws.Range('D27').Value2 = '25.45'
$ws.Range('E27').Value2 = '  -2.55%  '
$ws.Range('B28').Value2 = 'Cosmos'
$ws.Range('C28').Value2 = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D28').Value2 = '9.95'
$ws.Range('E28').Value2 = '  +1.66%  '
$ws.Range('B29').Value2 = 'Toncoin'
$ws.Range('C29').Value2 = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value2 = '2.18'
$ws.Range('E29').Value2 = '  -1.44%  '
$ws.Range('D30').Value2 = '0.139'
$ws.Range('E30').Value2 = '  +7.23%  '
$ws.Range('D31').Value2 = '34.45'
$ws.Range('E31').Value2 = '  -2.30%  '
$ws.Range('D32').Value2 = '49.39'
$ws.Range('E32').Value2 = '  -0.48%  '
$ws.Range('D33').Value2 = '19.77'
$ws.Range('E33').Value2 = '  +0.71%  '
$ws.Range('D34').Value2 = '5.30'
$ws.Range('E34').Value2 = '  -0.66%  '
$ws.Range('D35').Value2 = '0.0776'
$ws.Range('E35').Value2 = '  +1.02%  '
$ws.Range('D36').Value2 = '1.01'
$ws.Range('E36').Value2 = '  +0.32%  '
$ws.Range('D37').Value2 = '1.94'
$ws.Range('E37').Value2 = '  +1.63%  '
$ws.Range('E38').Value2 = '  -0.62%  '
$ws.Range('D39').Value2 = '2.91'
$ws.Range('E39').Value2 = '  -0.90%  '
$ws.Range('E40').Value2 = '  -0.44%  '
$ws.Range('D41').Value2 = '2.21'
$ws.Range('E41').Value2 = '  -0.66%  '
$ws.Range('D42').Value2 = '118.71'
$ws.Range('E42').Value2 = '  -3.99%  '
$ws.Range('D43').Value2 = '21.17'
$ws.Range('E43').Value2 = '  +2.42%  '
$ws.Range('D44').Value2 = '0.0291'
$ws.Range('E44').Value2 = '  -0.75%  '
$ws.Range('D45').Value2 = '1.969.40'
$ws.Range('E45').Value2 = '  -0.67%  '
$ws.Range('D46').Value2 = '2.97'
$ws.Range('E46').Value2 = '  -0.17%  '
$ws.Range('E47').Value2 = '  -3.94%  '
$ws.Range('D48').Value2 = '9.00'
$ws.Range('E48').Value2 = '  +0.54%  '
$ws.Range('E49').Value2 = '  -3.83%  '
$ws.Range('D50').Value2 = '5.09'
$ws.Range('E50').Value2 = '  -5.23%  '
$ws.Range('D51').Value2 = '56.33'
$ws.Range('E51').Value2 = '  +3.49%  '
